# Apply the "add all attribute level tables" edit to get_level_from_xp.xlsx
#
# Before: a single sheet "Feuil1" holding the get_level_from_xp table.
# After : four new lookup sheets (defense, max_health, speed, strength) are
#         inserted before it, and the original sheet is renamed to
#         "get_level_from_xp".

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet -------------------------------------------------
$main = $wb.Worksheets.Item(1)
$main.Name = "get_level_from_xp"

$uuid = "891dde02-2eec-45eb-8b4a-a68f3b0acc7a"

# Data for the four new sheets, in the order they must appear in the workbook
# (and therefore in sharedStrings.xml too).
$defs = @(
    @{ Name = "defense";    Score = "rpg_lvl_defense";    Attribute = "generic.armor";         Step = 1;     Decimal = $false },
    @{ Name = "max_health"; Score = "rpg_lvl_max_health"; Attribute = "generic.max_health";     Step = 1;     Decimal = $false },
    @{ Name = "speed";      Score = "rpg_lvl_speed";      Attribute = "generic.movement_speed"; Step = 0.003; Decimal = $true },
    @{ Name = "strength";   Score = "rpg_lvl_strength";   Attribute = "generic.attack_damage";  Step = 0.01;  Decimal = $true }
)

$prev = $main
foreach ($def in $defs) {
    $ws = $wb.Worksheets.Add($null, $prev)
    $ws.Name = $def.Name

    $ws.Range("A1").Value = $def.Score
    $ws.Range("B1").Value = $def.Attribute
    $ws.Range("C1").Value = $def.Step

    $ws.Range("B2").Value = 1
    $ws.Range("C2").Formula = "=C1"

    $ws.Range("B3").Value = 2
    $ws.Range("C3").Formula = "=C2+`$C`$1"

    $ws.Range("B4:B26").Formula = "=B3+1"
    $ws.Range("C4:C26").Formula = "=C3+`$C`$1"

    if ($def.Decimal) {
        $formula = '="execute if score @s " & $A$1 & " matches " & B2 & " run attribute @s " & $B$1 & " modifier add ' + $uuid + ' rpg_stats " & SUBSTITUTE(C2,",",".") & " add"'
        $formula3 = '="execute if score @s " & $A$1 & " matches " & B3 & " run attribute @s " & $B$1 & " modifier add ' + $uuid + ' rpg_stats " & SUBSTITUTE(C3,",",".") & " add"'
    } else {
        $formula = '="execute if score @s " & $A$1 & " matches " & B2 & " run attribute @s " & $B$1 & " modifier add ' + $uuid + ' rpg_stats " & C2 & " add"'
        $formula3 = '="execute if score @s " & $A$1 & " matches " & B3 & " run attribute @s " & $B$1 & " modifier add ' + $uuid + ' rpg_stats " & C3 & " add"'
    }

    $ws.Range("D2").Formula = $formula
    $ws.Range("D3:D26").Formula = $formula3

    $ws.Cells.Item(2, 4).Select()

    $ws.Columns("A:D").AutoFit() | Out-Null

    $prev = $ws
}

# the original sheet should end up last, after all four new lookup sheets
$main.Move($null, $prev)

$wb.Worksheets.Item("get_level_from_xp").Select()
